# ============================================================
# Edit script: restructure PlayerPerformance workbook
# - Add 'Player Info' sheet (first)
# - Rename MATCH_CARD_LINK -> MATCH_CODE on 'ODI Batting' and 'ODI Bowling',
#   converting the URL values to bare match-code numbers (stored as text)
# - Add 'ODI Batting Extra' sheet (last)
# ============================================================

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $text) {
    # Force a digit-looking string to be stored as TEXT, not coerced to a number,
    # then strip the resulting quote-prefix style so no stray format is left behind.
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

function Set-HeaderStyle($rng) {
    $rng.Font.Bold = $true
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
}

# --- 1) Insert 'Player Info' sheet before 'ODI Batting' ---
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet, $null)
$playerInfo.Name = "Player Info"

$playerInfo.Cells.Item(1,1).Value = "ID"
$playerInfo.Cells.Item(1,2).Value = "NAME"
$playerInfo.Cells.Item(1,3).Value = "BATTING_HAND"
$playerInfo.Cells.Item(1,4).Value = "BOWL_STYLE"
Set-HeaderStyle($playerInfo.Range("A1:D1"))

Set-TextValue $playerInfo.Cells.Item(2,1) "4423"
$playerInfo.Cells.Item(2,2).Value = "Imad Wasim"
$playerInfo.Cells.Item(2,3).Value = "Left Handed"
$playerInfo.Cells.Item(2,4).Value = "Left Arm Orthodox"

# --- 2) 'ODI Batting': MATCH_CARD_LINK header+values -> MATCH_CODE ---
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Cells.Item(1,4).Value = "MATCH_CODE"
Set-TextValue $battingSheet.Cells.Item(2,4) "3820"
Set-TextValue $battingSheet.Cells.Item(3,4) "3821"
Set-TextValue $battingSheet.Cells.Item(4,4) "3822"
Set-TextValue $battingSheet.Cells.Item(5,4) "3836"
Set-TextValue $battingSheet.Cells.Item(6,4) "3838"
Set-TextValue $battingSheet.Cells.Item(7,4) "3879"
Set-TextValue $battingSheet.Cells.Item(8,4) "3921"
Set-TextValue $battingSheet.Cells.Item(9,4) "3925"
Set-TextValue $battingSheet.Cells.Item(10,4) "3926"
Set-TextValue $battingSheet.Cells.Item(11,4) "3930"
Set-TextValue $battingSheet.Cells.Item(12,4) "3932"
Set-TextValue $battingSheet.Cells.Item(13,4) "3939"
Set-TextValue $battingSheet.Cells.Item(14,4) "3943"
Set-TextValue $battingSheet.Cells.Item(15,4) "3944"
Set-TextValue $battingSheet.Cells.Item(16,4) "3972"
Set-TextValue $battingSheet.Cells.Item(17,4) "3973"
Set-TextValue $battingSheet.Cells.Item(18,4) "3975"
Set-TextValue $battingSheet.Cells.Item(19,4) "3977"
Set-TextValue $battingSheet.Cells.Item(20,4) "4017"
Set-TextValue $battingSheet.Cells.Item(21,4) "4018"
Set-TextValue $battingSheet.Cells.Item(22,4) "4019"
Set-TextValue $battingSheet.Cells.Item(23,4) "4034"
Set-TextValue $battingSheet.Cells.Item(24,4) "4037"
Set-TextValue $battingSheet.Cells.Item(25,4) "4044"
Set-TextValue $battingSheet.Cells.Item(26,4) "4045"
Set-TextValue $battingSheet.Cells.Item(27,4) "4050"
Set-TextValue $battingSheet.Cells.Item(28,4) "4079"
Set-TextValue $battingSheet.Cells.Item(29,4) "4081"
Set-TextValue $battingSheet.Cells.Item(30,4) "4084"
Set-TextValue $battingSheet.Cells.Item(31,4) "4087"
Set-TextValue $battingSheet.Cells.Item(32,4) "4223"
Set-TextValue $battingSheet.Cells.Item(33,4) "4225"
Set-TextValue $battingSheet.Cells.Item(34,4) "4237"
Set-TextValue $battingSheet.Cells.Item(35,4) "4241"
Set-TextValue $battingSheet.Cells.Item(36,4) "4244"
Set-TextValue $battingSheet.Cells.Item(37,4) "4247"
Set-TextValue $battingSheet.Cells.Item(38,4) "4273"
Set-TextValue $battingSheet.Cells.Item(39,4) "4274"
Set-TextValue $battingSheet.Cells.Item(40,4) "4275"
Set-TextValue $battingSheet.Cells.Item(41,4) "4276"
Set-TextValue $battingSheet.Cells.Item(42,4) "4277"
Set-TextValue $battingSheet.Cells.Item(43,4) "4287"
Set-TextValue $battingSheet.Cells.Item(44,4) "4292"
Set-TextValue $battingSheet.Cells.Item(45,4) "4294"
Set-TextValue $battingSheet.Cells.Item(46,4) "4297"
Set-TextValue $battingSheet.Cells.Item(47,4) "4300"
Set-TextValue $battingSheet.Cells.Item(48,4) "4304"
Set-TextValue $battingSheet.Cells.Item(49,4) "4324"
Set-TextValue $battingSheet.Cells.Item(50,4) "4334"
Set-TextValue $battingSheet.Cells.Item(51,4) "4337"
Set-TextValue $battingSheet.Cells.Item(52,4) "4340"
Set-TextValue $battingSheet.Cells.Item(53,4) "4349"
Set-TextValue $battingSheet.Cells.Item(54,4) "4375"
Set-TextValue $battingSheet.Cells.Item(55,4) "4432"
Set-TextValue $battingSheet.Cells.Item(56,4) "4433"

# --- 3) 'ODI Bowling': MATCH_CARD_LINK header+values -> MATCH_CODE ---
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Cells.Item(1,2).Value = "MATCH_CODE"
Set-TextValue $bowlingSheet.Cells.Item(2,2) "3820"
Set-TextValue $bowlingSheet.Cells.Item(3,2) "3821"
Set-TextValue $bowlingSheet.Cells.Item(4,2) "3822"
Set-TextValue $bowlingSheet.Cells.Item(5,2) "3836"
Set-TextValue $bowlingSheet.Cells.Item(6,2) "3838"
Set-TextValue $bowlingSheet.Cells.Item(7,2) "3879"
Set-TextValue $bowlingSheet.Cells.Item(8,2) "3921"
Set-TextValue $bowlingSheet.Cells.Item(9,2) "3925"
Set-TextValue $bowlingSheet.Cells.Item(10,2) "3926"
Set-TextValue $bowlingSheet.Cells.Item(11,2) "3930"
Set-TextValue $bowlingSheet.Cells.Item(12,2) "3932"
Set-TextValue $bowlingSheet.Cells.Item(13,2) "3939"
Set-TextValue $bowlingSheet.Cells.Item(14,2) "3943"
Set-TextValue $bowlingSheet.Cells.Item(15,2) "3944"
Set-TextValue $bowlingSheet.Cells.Item(16,2) "3972"
Set-TextValue $bowlingSheet.Cells.Item(17,2) "3973"
Set-TextValue $bowlingSheet.Cells.Item(18,2) "3975"
Set-TextValue $bowlingSheet.Cells.Item(19,2) "3977"
Set-TextValue $bowlingSheet.Cells.Item(20,2) "4017"
Set-TextValue $bowlingSheet.Cells.Item(21,2) "4018"
Set-TextValue $bowlingSheet.Cells.Item(22,2) "4019"
Set-TextValue $bowlingSheet.Cells.Item(23,2) "4034"
Set-TextValue $bowlingSheet.Cells.Item(24,2) "4037"
Set-TextValue $bowlingSheet.Cells.Item(25,2) "4044"
Set-TextValue $bowlingSheet.Cells.Item(26,2) "4045"
Set-TextValue $bowlingSheet.Cells.Item(27,2) "4050"
Set-TextValue $bowlingSheet.Cells.Item(28,2) "4079"
Set-TextValue $bowlingSheet.Cells.Item(29,2) "4081"
Set-TextValue $bowlingSheet.Cells.Item(30,2) "4084"
Set-TextValue $bowlingSheet.Cells.Item(31,2) "4087"
Set-TextValue $bowlingSheet.Cells.Item(32,2) "4223"
Set-TextValue $bowlingSheet.Cells.Item(33,2) "4225"
Set-TextValue $bowlingSheet.Cells.Item(34,2) "4237"
Set-TextValue $bowlingSheet.Cells.Item(35,2) "4241"
Set-TextValue $bowlingSheet.Cells.Item(36,2) "4244"
Set-TextValue $bowlingSheet.Cells.Item(37,2) "4247"
Set-TextValue $bowlingSheet.Cells.Item(38,2) "4273"
Set-TextValue $bowlingSheet.Cells.Item(39,2) "4274"
Set-TextValue $bowlingSheet.Cells.Item(40,2) "4275"
Set-TextValue $bowlingSheet.Cells.Item(41,2) "4276"
Set-TextValue $bowlingSheet.Cells.Item(42,2) "4277"
Set-TextValue $bowlingSheet.Cells.Item(43,2) "4292"
Set-TextValue $bowlingSheet.Cells.Item(44,2) "4294"
Set-TextValue $bowlingSheet.Cells.Item(45,2) "4297"
Set-TextValue $bowlingSheet.Cells.Item(46,2) "4300"
Set-TextValue $bowlingSheet.Cells.Item(47,2) "4324"
Set-TextValue $bowlingSheet.Cells.Item(48,2) "4334"
Set-TextValue $bowlingSheet.Cells.Item(49,2) "4337"
Set-TextValue $bowlingSheet.Cells.Item(50,2) "4340"
Set-TextValue $bowlingSheet.Cells.Item(51,2) "4349"
Set-TextValue $bowlingSheet.Cells.Item(52,2) "4375"
Set-TextValue $bowlingSheet.Cells.Item(53,2) "4432"
Set-TextValue $bowlingSheet.Cells.Item(54,2) "4433"

# --- 4) Append 'ODI Batting Extra' sheet after 'ODI Bowling' ---
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$battingExtra = $wb.Worksheets.Add($null, $bowlingSheet)
$battingExtra.Name = "ODI Batting Extra"

$battingExtra.Cells.Item(1,1).Value = "MATCH_CODE"
$battingExtra.Cells.Item(1,2).Value = "BATTING_POSITION"
$battingExtra.Cells.Item(1,3).Value = "NUM_4"
$battingExtra.Cells.Item(1,4).Value = "NUM_6"
$battingExtra.Cells.Item(1,5).Value = "PERCENT_RUNS_OF_TOTAL"
$battingExtra.Cells.Item(1,6).Value = "MAN_OF_MATCH"
Set-HeaderStyle($battingExtra.Range("A1:F1"))

# row 2: match 4247
Set-TextValue $battingExtra.Cells.Item(2,1) "4247"
$battingExtra.Cells.Item(2,2).Value = ""
$battingExtra.Cells.Item(2,3).Value = ""
$battingExtra.Cells.Item(2,4).Value = ""
$battingExtra.Cells.Item(2,5).Value = ""
$battingExtra.Cells.Item(2,6).Value = "NO"

# row 3: match 4273
Set-TextValue $battingExtra.Cells.Item(3,1) "4273"
$battingExtra.Cells.Item(3,2).Value = 7
Set-TextValue $battingExtra.Cells.Item(3,3) "4"
Set-TextValue $battingExtra.Cells.Item(3,4) "1"
Set-TextValue $battingExtra.Cells.Item(3,5) "10.00%"
$battingExtra.Cells.Item(3,6).Value = "NO"

# row 4: match 4274
Set-TextValue $battingExtra.Cells.Item(4,1) "4274"
$battingExtra.Cells.Item(4,2).Value = 8
Set-TextValue $battingExtra.Cells.Item(4,3) "2"
Set-TextValue $battingExtra.Cells.Item(4,4) "1"
Set-TextValue $battingExtra.Cells.Item(4,5) "6.69%"
$battingExtra.Cells.Item(4,6).Value = "NO"

# row 5: match 4275
Set-TextValue $battingExtra.Cells.Item(5,1) "4275"
$battingExtra.Cells.Item(5,2).Value = 7
Set-TextValue $battingExtra.Cells.Item(5,3) "4"
Set-TextValue $battingExtra.Cells.Item(5,4) "0"
Set-TextValue $battingExtra.Cells.Item(5,5) "23.12%"
$battingExtra.Cells.Item(5,6).Value = "NO"

# row 6: match 4276
Set-TextValue $battingExtra.Cells.Item(6,1) "4276"
$battingExtra.Cells.Item(6,2).Value = ""
$battingExtra.Cells.Item(6,3).Value = ""
$battingExtra.Cells.Item(6,4).Value = ""
$battingExtra.Cells.Item(6,5).Value = ""
$battingExtra.Cells.Item(6,6).Value = "NO"

# row 7: match 4277
Set-TextValue $battingExtra.Cells.Item(7,1) "4277"
$battingExtra.Cells.Item(7,2).Value = 7
Set-TextValue $battingExtra.Cells.Item(7,3) "6"
Set-TextValue $battingExtra.Cells.Item(7,4) "1"
Set-TextValue $battingExtra.Cells.Item(7,5) "16.29%"
$battingExtra.Cells.Item(7,6).Value = "NO"

# row 8: match 4287
Set-TextValue $battingExtra.Cells.Item(8,1) "4287"
$battingExtra.Cells.Item(8,2).Value = ""
$battingExtra.Cells.Item(8,3).Value = ""
$battingExtra.Cells.Item(8,4).Value = ""
$battingExtra.Cells.Item(8,5).Value = ""
$battingExtra.Cells.Item(8,6).Value = "NO"

# row 9: match 4292
Set-TextValue $battingExtra.Cells.Item(9,1) "4292"
$battingExtra.Cells.Item(9,2).Value = 7
Set-TextValue $battingExtra.Cells.Item(9,3) "0"
Set-TextValue $battingExtra.Cells.Item(9,4) "1"
Set-TextValue $battingExtra.Cells.Item(9,5) "2.22%"
$battingExtra.Cells.Item(9,6).Value = "NO"

# row 10: match 4294
Set-TextValue $battingExtra.Cells.Item(10,1) "4294"
$battingExtra.Cells.Item(10,2).Value = ""
$battingExtra.Cells.Item(10,3).Value = ""
$battingExtra.Cells.Item(10,4).Value = ""
$battingExtra.Cells.Item(10,5).Value = ""
$battingExtra.Cells.Item(10,6).Value = "NO"

# row 11: match 4297
Set-TextValue $battingExtra.Cells.Item(11,1) "4297"
$battingExtra.Cells.Item(11,2).Value = 7
Set-TextValue $battingExtra.Cells.Item(11,3) "1"
Set-TextValue $battingExtra.Cells.Item(11,4) "1"
Set-TextValue $battingExtra.Cells.Item(11,5) "3.53%"
$battingExtra.Cells.Item(11,6).Value = "NO"

# row 12: match 4300
Set-TextValue $battingExtra.Cells.Item(12,1) "4300"
$battingExtra.Cells.Item(12,2).Value = 8
Set-TextValue $battingExtra.Cells.Item(12,3) "3"
Set-TextValue $battingExtra.Cells.Item(12,4) "0"
Set-TextValue $battingExtra.Cells.Item(12,5) "8.42%"
$battingExtra.Cells.Item(12,6).Value = "NO"

# row 13: match 4304
Set-TextValue $battingExtra.Cells.Item(13,1) "4304"
$battingExtra.Cells.Item(13,2).Value = 7
Set-TextValue $battingExtra.Cells.Item(13,3) "0"
Set-TextValue $battingExtra.Cells.Item(13,4) "0"
Set-TextValue $battingExtra.Cells.Item(13,5) "0.95%"
$battingExtra.Cells.Item(13,6).Value = "NO"

# row 14: match 4324
Set-TextValue $battingExtra.Cells.Item(14,1) "4324"
$battingExtra.Cells.Item(14,2).Value = 7
Set-TextValue $battingExtra.Cells.Item(14,3) "6"
Set-TextValue $battingExtra.Cells.Item(14,4) "0"
Set-TextValue $battingExtra.Cells.Item(14,5) "21.70%"
$battingExtra.Cells.Item(14,6).Value = "NO"

# row 15: match 4334
Set-TextValue $battingExtra.Cells.Item(15,1) "4334"
$battingExtra.Cells.Item(15,2).Value = 6
Set-TextValue $battingExtra.Cells.Item(15,3) "3"
Set-TextValue $battingExtra.Cells.Item(15,4) "0"
Set-TextValue $battingExtra.Cells.Item(15,5) "7.47%"
$battingExtra.Cells.Item(15,6).Value = "NO"

# row 16: match 4337
Set-TextValue $battingExtra.Cells.Item(16,1) "4337"
$battingExtra.Cells.Item(16,2).Value = ""
$battingExtra.Cells.Item(16,3).Value = ""
$battingExtra.Cells.Item(16,4).Value = ""
$battingExtra.Cells.Item(16,5).Value = ""
$battingExtra.Cells.Item(16,6).Value = "NO"

# row 17: match 4340
Set-TextValue $battingExtra.Cells.Item(17,1) "4340"
$battingExtra.Cells.Item(17,2).Value = 7
Set-TextValue $battingExtra.Cells.Item(17,3) "5"
Set-TextValue $battingExtra.Cells.Item(17,4) "0"
Set-TextValue $battingExtra.Cells.Item(17,5) "21.30%"
$battingExtra.Cells.Item(17,6).Value = "YES"

# row 18: match 4349
Set-TextValue $battingExtra.Cells.Item(18,1) "4349"
$battingExtra.Cells.Item(18,2).Value = 6
Set-TextValue $battingExtra.Cells.Item(18,3) "6"
Set-TextValue $battingExtra.Cells.Item(18,4) "1"
Set-TextValue $battingExtra.Cells.Item(18,5) "13.65%"
$battingExtra.Cells.Item(18,6).Value = "NO"

# row 19: match 4375
Set-TextValue $battingExtra.Cells.Item(19,1) "4375"
$battingExtra.Cells.Item(19,2).Value = 7
Set-TextValue $battingExtra.Cells.Item(19,3) "2"
Set-TextValue $battingExtra.Cells.Item(19,4) "0"
Set-TextValue $battingExtra.Cells.Item(19,5) "3.93%"
$battingExtra.Cells.Item(19,6).Value = "NO"

# row 20: match 4432
Set-TextValue $battingExtra.Cells.Item(20,1) "4432"
$battingExtra.Cells.Item(20,2).Value = ""
$battingExtra.Cells.Item(20,3).Value = ""
$battingExtra.Cells.Item(20,4).Value = ""
$battingExtra.Cells.Item(20,5).Value = ""
$battingExtra.Cells.Item(20,6).Value = "NO"

# row 21: match 4433
Set-TextValue $battingExtra.Cells.Item(21,1) "4433"
$battingExtra.Cells.Item(21,2).Value = ""
$battingExtra.Cells.Item(21,3).Value = ""
$battingExtra.Cells.Item(21,4).Value = ""
$battingExtra.Cells.Item(21,5).Value = ""
$battingExtra.Cells.Item(21,6).Value = "NO"

# --- 5) Restore original active sheet/selection (first sheet) ---
$null = $wb.Worksheets.Item("Player Info").Activate()
$null = $wb.Worksheets.Item("Player Info").Range("A1").Select()

